$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.667.95"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.328.56"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.53"
$ws.Range("E6").Value = "  +3.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.28"
$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +2.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.13"
$ws.Range("E12").Value = "  +2.70%  "

$ws.Range("E13").Value = "  +0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.681.68"
$ws.Range("E14").Value = "  +4.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.62"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  +7.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.323.07"
$ws.Range("E17").Value = "  +4.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.645.42"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  +5.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("E20").Value = "  +6.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.26"
$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.92"
$ws.Range("E22").Value = "  +3.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.38"
$ws.Range("E24").Value = "  +7.16%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  +1.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.42"
$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.50"
$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.67"
$ws.Range("E30").Value = "  +9.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.28"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.38"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0900"
$ws.Range("E33").Value = "  -1.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.47"
$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0363"
$ws.Range("E36").Value = "  +3.52%  "

$ws.Range("E37").Value = "  -3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  +3.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.37"
$ws.Range("E40").Value = "  +9.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.236"
$ws.Range("E41").Value = "  +10.48%  "

$ws.Range("E42").Value = "  +18.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.10"
$ws.Range("E43").Value = "  -4.11%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.17"
$ws.Range("E44").Value = "  +9.20%  "

$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.28"
$ws.Range("E45").Value = "  -1.46%  "

$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("E47").Value = "  +4.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.45"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.192"
$ws.Range("E50").Value = "  +18.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.560.69"
$ws.Range("E51").Value = "  +4.30%  "
